# Apply updated Survived predictions (feature engineering + CNN based solution)
# to the gender_submission sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gender_submission")

    $ws.Cells.Item(5, 2).Value = 1
    $ws.Cells.Item(8, 2).Value = 0
    $ws.Cells.Item(19, 2).Value = 1
    $ws.Cells.Item(29, 2).Value = 1
    $ws.Cells.Item(35, 2).Value = 0
    $ws.Cells.Item(37, 2).Value = 1
    $ws.Cells.Item(38, 2).Value = 0
    $ws.Cells.Item(39, 2).Value = 0
    $ws.Cells.Item(51, 2).Value = 0
    $ws.Cells.Item(66, 2).Value = 0
    $ws.Cells.Item(68, 2).Value = 0
    $ws.Cells.Item(77, 2).Value = 0
    $ws.Cells.Item(82, 2).Value = 1
    $ws.Cells.Item(89, 2).Value = 0
    $ws.Cells.Item(90, 2).Value = 0
    $ws.Cells.Item(100, 2).Value = 0
    $ws.Cells.Item(113, 2).Value = 0
    $ws.Cells.Item(115, 2).Value = 0
    $ws.Cells.Item(134, 2).Value = 0
    $ws.Cells.Item(140, 2).Value = 0
    $ws.Cells.Item(144, 2).Value = 0
    $ws.Cells.Item(146, 2).Value = 1
    $ws.Cells.Item(159, 2).Value = 0
    $ws.Cells.Item(160, 2).Value = 1
    $ws.Cells.Item(161, 2).Value = 0
    $ws.Cells.Item(162, 2).Value = 0
    $ws.Cells.Item(163, 2).Value = 1
    $ws.Cells.Item(171, 2).Value = 0
    $ws.Cells.Item(198, 2).Value = 0
    $ws.Cells.Item(199, 2).Value = 0
    $ws.Cells.Item(201, 2).Value = 0
    $ws.Cells.Item(204, 2).Value = 1
    $ws.Cells.Item(215, 2).Value = 0
    $ws.Cells.Item(216, 2).Value = 1
    $ws.Cells.Item(219, 2).Value = 0
    $ws.Cells.Item(225, 2).Value = 1
    $ws.Cells.Item(227, 2).Value = 0
    $ws.Cells.Item(229, 2).Value = 0
    $ws.Cells.Item(238, 2).Value = 1
    $ws.Cells.Item(239, 2).Value = 1
    $ws.Cells.Item(251, 2).Value = 0
    $ws.Cells.Item(254, 2).Value = 0
    $ws.Cells.Item(260, 2).Value = 0
    $ws.Cells.Item(270, 2).Value = 0
    $ws.Cells.Item(282, 2).Value = 0
    $ws.Cells.Item(284, 2).Value = 0
    $ws.Cells.Item(285, 2).Value = 0
    $ws.Cells.Item(286, 2).Value = 0
    $ws.Cells.Item(293, 2).Value = 0
    $ws.Cells.Item(295, 2).Value = 1
    $ws.Cells.Item(306, 2).Value = 0
    $ws.Cells.Item(313, 2).Value = 1
    $ws.Cells.Item(318, 2).Value = 1
    $ws.Cells.Item(319, 2).Value = 1
    $ws.Cells.Item(365, 2).Value = 1
    $ws.Cells.Item(369, 2).Value = 0
    $ws.Cells.Item(371, 2).Value = 1
    $ws.Cells.Item(378, 2).Value = 0
    $ws.Cells.Item(385, 2).Value = 0
    $ws.Cells.Item(392, 2).Value = 0
    $ws.Cells.Item(405, 2).Value = 0
    $ws.Cells.Item(409, 2).Value = 0
    $ws.Cells.Item(410, 2).Value = 0
    $ws.Cells.Item(411, 2).Value = 0
    $ws.Cells.Item(412, 2).Value = 0
